$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.375.04"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.190.05"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.75"
$ws.Range("E5").Value = "  +5.52%  "
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +4.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.64"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.08"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0938"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("E13").Value = "  +6.11%  "
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.516.40"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.173.26"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.307.72"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0955"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.08"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.90"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.94"
$ws.Range("E25").Value = "  +9.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.83"
$ws.Range("E26").Value = "  +21.19%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("E29").Value = "  -4.27%  "
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.12"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0758"
$ws.Range("E34").Value = "  +7.03%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.49"
$ws.Range("E36").Value = "  +5.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.71"
$ws.Range("E37").Value = "  +12.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.18"
$ws.Range("E38").Value = "  +7.41%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0298"
$ws.Range("E40").Value = "  +7.24%  "
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("E42").Value = "  +16.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.69"
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.21"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.05"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.201"
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.67"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("E49").Value = "  +9.60%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +0.78%  "
